$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.432.83'
$ws.Range("E2").Value = '  +1.26%  '

# Row 3
$ws.Range("D3").Value = '3.156.70'
$ws.Range("E3").Value = '  +1.38%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.89'
$ws.Range("E5").Value = '  +1.49%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.19'
$ws.Range("E6").Value = '  +1.17%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("D8").Value = '3.146.51'
$ws.Range("E8").Value = '  +1.30%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  +0.24%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("E10").Value = '  +1.90%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.95'
$ws.Range("E11").Value = '  +4.57%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.461'
$ws.Range("E12").Value = '  +0.10%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  -0.14%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.45'
$ws.Range("E14").Value = '  -0.25%  '

# Row 15
$ws.Range("D15").Value = '3.678.53'
$ws.Range("E15").Value = '  +1.37%  '

# Row 16
$ws.Range("E16").Value = '  -0.51%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.30'
$ws.Range("E17").Value = '  +2.35%  '

# Row 18
$ws.Range("D18").Value = '64.165.03'
$ws.Range("E18").Value = '  +0.93%  '

# Row 19
$ws.Range("D19").Value = '3.150.85'
$ws.Range("E19").Value = '  +1.44%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '470.48'
$ws.Range("E20").Value = '  +1.58%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.42'
$ws.Range("E21").Value = '  +0.43%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.738'
$ws.Range("E22").Value = '  +1.11%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.62'
$ws.Range("E23").Value = '  +0.99%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.42'
$ws.Range("E24").Value = '  +12.40%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.16'
$ws.Range("E25").Value = '  +0.32%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.39'
$ws.Range("E26").Value = '  -0.07%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.95'
$ws.Range("E28").Value = '  +10.76%  '

# Row 29
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.27'
$ws.Range("E29").Value = '  +1.93%  '

# Row 30
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.45'
$ws.Range("E30").Value = '  +9.08%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.73'
$ws.Range("E31").Value = '  +1.78%  '

# Row 32
$ws.Range("E32").Value = '  +0.09%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.75'
$ws.Range("E33").Value = '  +3.16%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.113'
$ws.Range("E34").Value = '  +4.57%  '

# Row 35
$ws.Range("D35").Value = '0.0₃0872'
$ws.Range("E35").Value = '  +1.60%  '

# Row 36
$ws.Range("E36").Value = '  +2.52%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.20'
$ws.Range("E37").Value = '  +2.50%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.31'
$ws.Range("E38").Value = '  -1.09%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.24'
$ws.Range("E39").Value = '  -4.26%  '

# Row 40
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '462.58'
$ws.Range("E40").Value = '  +4.48%  '

# Row 41
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.44'
$ws.Range("E41").Value = '  +1.16%  '

# Row 42
$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.37'
$ws.Range("E42").Value = '  +7.56%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.300'
$ws.Range("E43").Value = '  +9.49%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0377'
$ws.Range("E44").Value = '  +1.86%  '

# Row 45
$ws.Range("D45").Value = '2.921.55'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.36'
$ws.Range("E46").Value = '  +12.74%  '

# Row 47
$ws.Range("E47").Value = '  -0.63%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.77'
$ws.Range("E48").Value = '  +8.64%  '

# Row 49
$ws.Range("E49").Value = '  +0.02%  '

# Row 50
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.26'
$ws.Range("E50").Value = '  +4.15%  '

# Row 51
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.111'
$ws.Range("E51").Value = '  +0.45%  '
